# ---------------------------------------------------------------------------
# KHL referee stats refresh (2025-12-18 23:30:06 snapshot)
#
# For both the "Главные" (referees) and "Линейные" (linesmen) sheets:
#   - officials who worked an additional game since the last publish get
#     their Games_KHL / PIM_* counters bumped to the latest totals
#   - every data row (2-26) gets its as_of_utc refresh timestamp updated
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$asOfUtc = "2025-12-18 23:30:06"
$firstDataRow = 2
$lastDataRow = 26

# --- Updated counters for "Главные" ---
$rowUpdatesMain = @(
    @{ Row = 4; Cells = @{ "C" = 26; "D" = 413; "E" = 179; "F" = 234; "G" = 15.88; "H" = 6.88; "I" = 9; "J" = 87; "K" = 102 } }
    @{ Row = 6; Cells = @{ "C" = 33; "D" = 561; "E" = 268; "F" = 293; "G" = 17; "H" = 8.119999999999999; "I" = 8.880000000000001; "J" = 114; "K" = 124 } }
    @{ Row = 8; Cells = @{ "C" = 31; "D" = 621; "E" = 324; "F" = 297; "G" = 20.03; "H" = 10.45; "I" = 9.58; "J" = 127; "K" = 126; "P" = 2; "W" = 14 } }
    @{ Row = 13; Cells = @{ "C" = 13; "D" = 197; "E" = 110; "F" = 87; "G" = 15.15; "H" = 8.460000000000001; "I" = 6.69; "J" = 55; "K" = 41 } }
    @{ Row = 16; Cells = @{ "C" = 32; "D" = 577; "E" = 288; "F" = 289; "G" = 18.03; "H" = 9; "I" = 9.029999999999999; "J" = 114; "K" = 112 } }
    @{ Row = 17; Cells = @{ "C" = 22; "D" = 344; "E" = 132; "F" = 212; "G" = 15.64; "H" = 6; "I" = 9.640000000000001; "J" = 61; "K" = 86 } }
    @{ Row = 20; Cells = @{ "C" = 33; "D" = 555; "E" = 246; "F" = 309; "G" = 16.82; "H" = 7.45; "I" = 9.359999999999999; "J" = 113; "K" = 117; "V" = 18; "W" = 12 } }
    @{ Row = 21; Cells = @{ "C" = 31; "D" = 452; "E" = 204; "F" = 248; "G" = 14.58; "H" = 6.58; "I" = 8; "J" = 92; "K" = 109; "V" = 10; "W" = 10 } }
    @{ Row = 25; Cells = @{ "C" = 34; "D" = 527; "E" = 260; "F" = 267; "G" = 15.5; "H" = 7.65; "I" = 7.85; "J" = 125; "K" = 126; "V" = 14; "W" = 18 } }
    @{ Row = 26; Cells = @{ "C" = 15; "D" = 318; "E" = 157; "F" = 161; "G" = 21.2; "H" = 10.47; "I" = 10.73; "J" = 61; "K" = 58; "P" = 1; "W" = 4 } }
)

# --- Updated counters for "Линейные" ---
$rowUpdatesLinear = @(
    @{ Row = 2; Cells = @{ "C" = 22; "D" = 440; "E" = 211; "F" = 229; "G" = 20; "H" = 9.59; "I" = 10.41; "J" = 83; "K" = 87; "P" = 1; "W" = 14 } }
    @{ Row = 3; Cells = @{ "C" = 34; "D" = 572; "E" = 290; "F" = 282; "G" = 16.82; "H" = 8.529999999999999; "I" = 8.289999999999999; "J" = 125; "K" = 111; "V" = 18; "W" = 14 } }
    @{ Row = 14; Cells = @{ "C" = 31; "D" = 498; "E" = 257; "F" = 241; "G" = 16.06; "H" = 8.289999999999999; "I" = 7.77; "J" = 126; "K" = 113; "V" = 16; "W" = 6 } }
    @{ Row = 16; Cells = @{ "C" = 32; "D" = 508; "E" = 238; "F" = 270; "G" = 15.88; "H" = 7.44; "I" = 8.44; "J" = 109; "K" = 115 } }
    @{ Row = 19; Cells = @{ "C" = 31; "D" = 489; "E" = 238; "F" = 251; "G" = 15.77; "H" = 7.68; "I" = 8.1; "J" = 114; "K" = 113 } }
    @{ Row = 20; Cells = @{ "C" = 21; "D" = 342; "E" = 191; "F" = 151; "G" = 16.29; "H" = 9.1; "I" = 7.19; "J" = 83; "K" = 73; "P" = 1; "W" = 8 } }
)

$sheetUpdates = @(
    @{ Name = "Главные"; Updates = $rowUpdatesMain }
    @{ Name = "Линейные"; Updates = $rowUpdatesLinear }
)

foreach ($sheetDef in $sheetUpdates) {
    $ws = $wb.Worksheets.Item($sheetDef.Name)

    # Apply the changed counters for the officials with new games played.
    foreach ($rowDef in $sheetDef.Updates) {
        foreach ($col in $rowDef.Cells.Keys) {
            $ws.Range("$col$($rowDef.Row)").Value = $rowDef.Cells[$col]
        }
    }

    # Stamp the refresh time on every data row.
    for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
        $ws.Range("AA$row").Value = $asOfUtc
    }
}
